$d = $word.ActiveDocument
$d.Content.Find.Execute("Saintization", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sanitization", 2)
